# comment from python script
#
# Fills in the "S4" evaluation column for each roster sheet (the next
# un-scored skill column), then leaves the selection/active sheet the way
# the grader last left them: Astronauta -> E22, Senador -> F20,
# Mago -> D8, Ninja active (selection stays at E4).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Astronauta")
$ws2 = $wb.Worksheets.Item("Senador")
$ws3 = $wb.Worksheets.Item("Mago")
$ws4 = $wb.Worksheets.Item("Ninja")

# --- Astronauta: new column E scores ---
$ws1.Range("E2").Value = 0
$ws1.Range("E3").Value = 0
$ws1.Range("E4").Value = 1
$ws1.Range("E5").Value = 1
$ws1.Range("E6").Value = 0.8
$ws1.Range("E8").Value = 1
$ws1.Range("E9").Value = 1
$ws1.Range("E10").Value = 0
$ws1.Range("E11").Value = 0
$ws1.Range("E12").Value = 0.5
$ws1.Range("E13").Value = 1
$ws1.Range("E14").Value = 1
$ws1.Range("E15").Value = 1
$ws1.Range("E16").Value = 1
$ws1.Range("E17").Value = 1
$ws1.Range("E18").Value = 1
$ws1.Range("E19").Value = 1
$ws1.Range("E20").Value = 1
$ws1.Range("E21").Value = 1

# --- Senador: new column E scores ---
$ws2.Range("E2").Value = 0
$ws2.Range("E3").Value = 1
$ws2.Range("E4").Value = 1
$ws2.Range("E5").Value = 1
$ws2.Range("E6").Value = 0.8
$ws2.Range("E8").Value = 1
$ws2.Range("E9").Value = 1
$ws2.Range("E10").Value = 0
$ws2.Range("E11").Value = 0
$ws2.Range("E12").Value = 0.5
$ws2.Range("E13").Value = 0.8
$ws2.Range("E14").Value = 1
$ws2.Range("E15").Value = 1
$ws2.Range("E16").Value = 1
$ws2.Range("E17").Value = 1
$ws2.Range("E18").Value = 0.7
$ws2.Range("E19").Value = 1
$ws2.Range("E20").Value = 0.9
$ws2.Range("E21").Value = 1

# --- Mago: new column D scores ---
$ws3.Range("D2").Value = 0.2
$ws3.Range("D3").Value = 1
$ws3.Range("D4").Value = 1
$ws3.Range("D5").Value = 1
$ws3.Range("D6").Value = 0.6
$ws3.Range("D8").Value = 1
$ws3.Range("D9").Value = 1
$ws3.Range("D10").Value = 0
$ws3.Range("D11").Value = 0
$ws3.Range("D12").Value = 0.4
$ws3.Range("D13").Value = 0.8
$ws3.Range("D14").Value = 1
$ws3.Range("D15").Value = 0.8
$ws3.Range("D16").Value = 1
$ws3.Range("D17").Value = 1
$ws3.Range("D18").Value = 0.8
$ws3.Range("D19").Value = 1
$ws3.Range("D20").Value = 0.5
$ws3.Range("D21").Value = 1

# --- Leave the same selection/active-sheet trail the author ended on ---
$ws1.Range("E22").Select() | Out-Null
$ws2.Range("F20").Select() | Out-Null
$ws3.Range("D8").Select() | Out-Null
$ws4.Activate() | Out-Null
